# [MIG] 12.0 account_bank_statement_import_adyen, account_bank_statement_clearing_account
#
# The migrated Adyen test fixture re-saves the same "Data" sheet with one
# substantive data correction: every Gross/Net Currency cell that used to
# read "EUR" is updated to "USD" (columns K and O, rows 5-44 - 77 cells in
# total). The GBP outlier (K28) is left untouched.
#
# The accompanying view-state tweaks (active selection moving to P38, the
# window's tab-ratio split) are also replayed here via the Excel object
# model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction -------------------------------------------------
# Replace the currency code across the whole sheet (only the Gross/Net
# Currency columns ever held "EUR", so a sheet-wide whole-cell replace is
# equivalent to, and less error-prone than, touching each cell by hand).
[void]$ws.Cells.Replace("EUR", "USD", 1)

# --- Window / view state ---------------------------------------------
$win = $excel.ActiveWindow
$win.TabRatio = 0.5

[void]$ws.Range("P38").Select()
